# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit-tracking tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2008
$ws.Range("I40").Value = 1265.75
$ws.Range("J40").Value = 2997.6667
$ws.Range("K40").Value = 1265.75
$ws.Range("L40").Value = 2997.6667
$ws.Range("M40").Value = -1090.75
$ws.Range("N40").Value = -3347.6667

$ws.Range("H53").Value = 210.15
$ws.Range("I53").Value = 181.2
$ws.Range("K53").Value = 181.2
$ws.Range("M53").Value = 455.8

$ws.Range("H96").Value = 1868
$ws.Range("J96").Value = 2433
$ws.Range("L96").Value = 7299
$ws.Range("N96").Value = -10045

$ws.Range("H100").Value = 1405.2307
$ws.Range("J100").Value = 728.75
$ws.Range("L100").Value = 728.75
$ws.Range("N100").Value = -1810.75

$ws.Range("H116").Value = 8800.941999999999
$ws.Range("I116").Value = 10245.4
$ws.Range("J116").Value = 6875
$ws.Range("K116").Value = 10245.4
$ws.Range("L116").Value = 6875
$ws.Range("M116").Value = -6803.4
$ws.Range("N116").Value = -13759

$ws.Range("H131").Value = 13709.533
$ws.Range("J131").Value = 28600.166
$ws.Range("L131").Value = 85800.49800000001
$ws.Range("N131").Value = -95880.49800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2997.7231
$ws.Range("I32").Value = 2903.9375
$ws.Range("K32").Value = 2903.9375
$ws.Range("M32").Value = -2616.9375

$ws.Range("H45").Value = 26079.6
$ws.Range("I45").Value = 37799.332
$ws.Range("K45").Value = 37799.332
$ws.Range("M45").Value = -37422.332

$ws.Range("H132").Value = 5361.3237
$ws.Range("I132").Value = 5017.0713
$ws.Range("K132").Value = 15051.2139
$ws.Range("M132").Value = -12521.2139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2579.8
$ws.Range("I105").Value = 3149.5
$ws.Range("K105").Value = 3149.5
$ws.Range("M105").Value = -1402.5

$ws.Range("H107").Value = 5840.8184
$ws.Range("I107").Value = 1504.6111
$ws.Range("J107").Value = 25353.75
$ws.Range("K107").Value = 1504.6111
$ws.Range("L107").Value = 25353.75
$ws.Range("M107").Value = 415.3888999999999
$ws.Range("N107").Value = -29193.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1278.8667
$ws.Range("I16").Value = 1227.4286
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 1227.4286
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -940.4286
$ws.Range("N16").Value = -2573

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H74").Value = 25299.8
$ws.Range("J74").Value = 37666
$ws.Range("L74").Value = 37666
$ws.Range("N74").Value = -39414

$ws.Range("H77").Value = 25299.8
$ws.Range("J77").Value = 37666
$ws.Range("L77").Value = 112998
$ws.Range("N77").Value = -121734

$ws.Range("H113").Value = 1278.8667
$ws.Range("I113").Value = 1227.4286
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 1227.4286
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 942.5714
$ws.Range("N113").Value = -6339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 500
$ws.Range("I48").Value = 500
$ws.Range("K48").Value = 1500
$ws.Range("M48").Value = -1250

$ws.Range("H86").Value = 569.8333
$ws.Range("I86").Value = 279.75
$ws.Range("J86").Value = 1150
$ws.Range("K86").Value = 839.25
$ws.Range("L86").Value = 3450
$ws.Range("M86").Value = 346.75
$ws.Range("N86").Value = -5822

$ws.Range("H89").Value = 569.8333
$ws.Range("I89").Value = 279.75
$ws.Range("J89").Value = 1150
$ws.Range("K89").Value = 2517.75
$ws.Range("L89").Value = 10350
$ws.Range("M89").Value = 3410.25
$ws.Range("N89").Value = -22206

$ws.Range("H109").Value = 1190.75
$ws.Range("I109").Value = 587.6667
$ws.Range("K109").Value = 1763.0001
$ws.Range("M109").Value = -723.0001

$ws.Range("H122").Value = 671.3333

$ws.Range("H129").Value = 1419.4
$ws.Range("I129").Value = 1024.25
$ws.Range("K129").Value = 3072.75
$ws.Range("M129").Value = 1927.25

$ws.Range("H131").Value = 57302.668
$ws.Range("I131").Value = 500500
$ws.Range("K131").Value = 1501500
$ws.Range("M131").Value = -1496460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9156.4
$ws.Range("J92").Value = 8945.5
$ws.Range("L92").Value = 8945.5
$ws.Range("N92").Value = -12689.5

$ws.Range("H107").Value = 702.44446
$ws.Range("I107").Value = 534.1429000000001
$ws.Range("K107").Value = 534.1429000000001
$ws.Range("M107").Value = 1385.8571

$ws.Range("H132").Value = 6274.3335
$ws.Range("I132").Value = 6476.8237
$ws.Range("K132").Value = 19430.4711
$ws.Range("M132").Value = -16900.4711

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1272.579
$ws.Range("I22").Value = 1105.6
$ws.Range("K22").Value = 1105.6
$ws.Range("M22").Value = -810.5999999999999

$ws.Range("H27").Value = 1272.579
$ws.Range("I27").Value = 1105.6
$ws.Range("K27").Value = 1105.6
$ws.Range("M27").Value = -998.5999999999999

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H108").Value = 32333.334
$ws.Range("J108").Value = 32333.334
$ws.Range("L108").Value = 32333.334
$ws.Range("N108").Value = -40013.334

$ws.Range("H132").Value = 5128.7896
$ws.Range("I132").Value = 4642.7144
$ws.Range("J132").Value = 6489.8
$ws.Range("K132").Value = 13928.1432
$ws.Range("L132").Value = 19469.4
$ws.Range("M132").Value = -11398.1432
$ws.Range("N132").Value = -24529.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9000.421
$ws.Range("I81").Value = 17268.166
$ws.Range("J81").Value = 5184.5386
$ws.Range("K81").Value = 34536.332
$ws.Range("L81").Value = 10369.0772
$ws.Range("M81").Value = -33475.332
$ws.Range("N81").Value = -12491.0772

$ws.Range("H84").Value = 9000.421
$ws.Range("I84").Value = 17268.166
$ws.Range("J84").Value = 5184.5386
$ws.Range("K84").Value = 172681.66
$ws.Range("L84").Value = 51845.386
$ws.Range("M84").Value = -167377.66
$ws.Range("N84").Value = -62453.386

$ws.Range("H107").Value = 276.15384
$ws.Range("I107").Value = 314.5
$ws.Range("J107").Value = 214.8
$ws.Range("K107").Value = 943.5
$ws.Range("L107").Value = 644.4000000000001
$ws.Range("M107").Value = 976.5
$ws.Range("N107").Value = -4484.4

$ws.Range("H122").Value = 2870.9092
$ws.Range("I122").Value = 2508.889
$ws.Range("K122").Value = 7526.667
$ws.Range("M122").Value = -5076.667
